$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Balance Sheet")
$cf = $wb.Worksheets.Item("Cash Flow")

$fBefore = $cf.Range("B15").Formula
Write-Host "Cash Flow B15 formula BEFORE insert:" $fBefore

$ws.Rows("51:51").Insert()

$fAfter = $cf.Range("B15").Formula
Write-Host "Cash Flow B15 formula AFTER insert:" $fAfter
